$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.778.25"
$ws.Range("E2").Value = "  +4.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.269.46"
$ws.Range("E3").Value = "  +3.56%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.13"
$ws.Range("E5").Value = "  +4.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.08"
$ws.Range("E6").Value = "  +7.60%  "
$ws.Range("E7").Value = "  +3.02%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  +5.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "54.45"
$ws.Range("E10").Value = "  +9.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.42"
$ws.Range("E11").Value = "  +8.34%  "
$ws.Range("E12").Value = "  +3.05%  "
$ws.Range("E13").Value = "  +3.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.67"
$ws.Range("E14").Value = "  +4.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.620.35"
$ws.Range("E15").Value = "  +3.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.19"
$ws.Range("E16").Value = "  +4.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.308.75"
$ws.Range("E17").Value = "  +5.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.755"
$ws.Range("E18").Value = "  +4.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.679.60"
$ws.Range("E19").Value = "  +4.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.38"
$ws.Range("E20").Value = "  +10.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("E21").Value = "  +3.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.93"
$ws.Range("E22").Value = "  +4.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.21"
$ws.Range("E23").Value = "  +3.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.78"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  +5.98%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +4.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.85"
$ws.Range("E28").Value = "  +3.68%  "
$ws.Range("E29").Value = "  +7.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.71"
$ws.Range("E30").Value = "  +6.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.18"
$ws.Range("E31").Value = "  +9.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.47"
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.20"
$ws.Range("E34").Value = "  +6.26%  "
$ws.Range("E35").Value = "  +5.09%  "
$ws.Range("E36").Value = "  +9.44%  "
$ws.Range("E38").Value = "  +10.82%  "
$ws.Range("E39").Value = "  +7.31%  "
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("E41").Value = "  +7.69%  "
$ws.Range("E42").Value = "  +7.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.26"
$ws.Range("E43").Value = "  +19.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.060.93"
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("E45").Value = "  +4.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.97"
$ws.Range("E46").Value = "  +12.34%  "
$ws.Range("E47").Value = "  +4.29%  "
$ws.Range("E48").Value = "  -4.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.492.28"
$ws.Range("E49").Value = "  +3.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.52"
$ws.Range("E50").Value = "  +4.52%  "
$ws.Range("E51").Value = "  +5.06%  "
